$d = $word.ActiveDocument

# --- Change 1: remove the existing "_GoBack" bookmark (after " con ") ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Change 2: change the "5,18" run's text to "4", keeping the run split
#     (so it doesn't get silently coalesced with its same-formatted
#     neighbours), then drop a fresh "_GoBack" bookmark right after it ---
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute("5,18", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $startPos = $findRange.Start
    $endPos = $findRange.End

    # Temporary bookmarks on both sides of the run keep it from merging
    # into its neighbouring runs (which share identical rPr) once its
    # text is rewritten.
    $d.Bookmarks.Add("ZZZ_TempLeft", $d.Range($startPos, $startPos)) | Out-Null
    $d.Bookmarks.Add("ZZZ_TempRight", $d.Range($endPos, $endPos)) | Out-Null

    $target = $d.Range($startPos, $endPos)
    $target.Text = "4"

    $d.Bookmarks.Item("ZZZ_TempLeft").Delete()
    $d.Bookmarks.Item("ZZZ_TempRight").Delete()

    # New end-of-run position (text shrank from "5,18" (4 chars) to "4" (1 char))
    $newEnd = $startPos + 1
    $d.Bookmarks.Add("_GoBack", $d.Range($newEnd, $newEnd)) | Out-Null
}
